# Updates the cryptos list values (prices, volume % changes, and a few
# coin name/link/price swaps) on sheet "Sheet1" of the workbook, mirroring
# the scheduled GitHub Actions data refresh described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that frequently looks numeric (e.g. "579.40",
# "0.0000277", "3.596.53"). Excel's COM layer auto-converts such strings to
# real numbers when assigned directly, which would silently drop trailing
# zeros or flip tiny values into scientific notation. Forcing the column to
# Text format before writing keeps every value as the exact original string,
# and resetting the style back to Normal afterwards avoids leaving a stray
# number-format override on the cells.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "69.826.79"
$ws.Range("E2").Value = "  -2.56%  "
$ws.Range("D3").Value = "3.564.77"
$ws.Range("E3").Value = "  -2.10%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "579.40"
$ws.Range("E5").Value = "  -1.85%  "
$ws.Range("D6").Value = "172.35"
$ws.Range("E6").Value = "  -4.64%  "
$ws.Range("D7").Value = "0.623"
$ws.Range("E7").Value = "  +1.59%  "
$ws.Range("D8").Value = "3.555.16"
$ws.Range("E8").Value = "  -2.08%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "0.191"
$ws.Range("E10").Value = "  -5.82%  "
$ws.Range("D11").Value = "6.51"
$ws.Range("E11").Value = "  +10.98%  "
$ws.Range("D12").Value = "0.605"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").Value = "47.63"
$ws.Range("E13").Value = "  -4.34%  "
$ws.Range("D14").Value = "0.0000277"
$ws.Range("E14").Value = "  -2.88%  "
$ws.Range("D15").Value = "691.87"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").Value = "4.160.15"
$ws.Range("E16").Value = "  -1.36%  "
$ws.Range("D17").Value = "8.86"
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.593.97"
$ws.Range("E18").Value = "  -1.88%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "69.947.67"
$ws.Range("E19").Value = "  -2.47%  "
$ws.Range("E20").Value = "  -0.91%  "
$ws.Range("D21").Value = "17.48"
$ws.Range("E21").Value = "  -4.52%  "
$ws.Range("D22").Value = "11.24"
$ws.Range("E22").Value = "  -3.32%  "
$ws.Range("D23").Value = "0.919"
$ws.Range("E23").Value = "  -1.75%  "
$ws.Range("D24").Value = "16.76"
$ws.Range("E24").Value = "  -5.75%  "
$ws.Range("D25").Value = "98.29"
$ws.Range("E25").Value = "  -4.88%  "
$ws.Range("D26").Value = "3.86"
$ws.Range("E26").Value = "  -4.04%  "
$ws.Range("D27").Value = "2.71"
$ws.Range("E27").Value = "  -4.74%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").Value = "9.51"
$ws.Range("E29").Value = "  -4.97%  "
$ws.Range("D30").Value = "33.89"
$ws.Range("E30").Value = "  -3.49%  "
$ws.Range("D31").Value = "8.98"
$ws.Range("E31").Value = "  -2.47%  "
$ws.Range("D32").Value = "3.20"
$ws.Range("E32").Value = "  -6.01%  "
$ws.Range("D33").Value = "7.35"
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("D34").Value = "1.35"
$ws.Range("E34").Value = "  -5.35%  "
$ws.Range("D35").Value = "3.86"
$ws.Range("E35").Value = "  -7.18%  "
$ws.Range("D36").Value = "573.24"
$ws.Range("E36").Value = "  -0.62%  "
$ws.Range("D37").Value = "10.91"
$ws.Range("E37").Value = "  -3.86%  "
$ws.Range("D38").Value = "0.106"
$ws.Range("E38").Value = "  -3.54%  "
$ws.Range("D39").Value = "57.89"
$ws.Range("E39").Value = "  -2.60%  "
$ws.Range("E40").Value = "  +0.36%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "0.139"
$ws.Range("E41").Value = "  -2.45%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.498.58"
$ws.Range("E42").Value = "  -4.91%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.0444"
$ws.Range("E43").Value = "  -6.06%  "
$ws.Range("D44").Value = "0.340"
$ws.Range("E44").Value = "  -2.26%  "
$ws.Range("D45").Value = "33.68"
$ws.Range("E45").Value = "  -5.75%  "
$ws.Range("D46").Value = "0.0₃0712"
$ws.Range("E46").Value = "  -6.33%  "
$ws.Range("D47").Value = "2.92"
$ws.Range("E47").Value = "  +3.22%  "
$ws.Range("D48").Value = "2.60"
$ws.Range("E48").Value = "  -5.81%  "
$ws.Range("D49").Value = "0.134"
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("D50").Value = "134.17"
$ws.Range("E50").Value = "  +1.85%  "
$ws.Range("D51").Value = "0.150"
$ws.Range("E51").Value = "  +0.33%  "

# Restore the default cell style now that the text has been written.
$priceRange.Style = "Normal"

